$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 1919.5333  # ALC!H40
$ws.Cells.Item(40, 9).Value = 1783.1666  # ALC!I40
$ws.Cells.Item(40, 10).Value = 2010.4445  # ALC!J40
$ws.Cells.Item(40, 11).Value = 1783.1666  # ALC!K40
$ws.Cells.Item(40, 12).Value = 2010.4445  # ALC!L40
$ws.Cells.Item(40, 13).Value = -1608.1666  # ALC!M40
$ws.Cells.Item(40, 14).Value = -2360.4445  # ALC!N40
$ws.Cells.Item(86, 8).Value = 3007.2  # ALC!H86
$ws.Cells.Item(86, 9).Value = 2844.3333  # ALC!I86
$ws.Cells.Item(86, 10).Value = 3495.8  # ALC!J86
$ws.Cells.Item(86, 11).Value = 2844.3333  # ALC!K86
$ws.Cells.Item(86, 12).Value = 3495.8  # ALC!L86
$ws.Cells.Item(86, 13).Value = -1721.3333  # ALC!M86
$ws.Cells.Item(86, 14).Value = -5741.8  # ALC!N86
$ws.Cells.Item(89, 8).Value = 3007.2  # ALC!H89
$ws.Cells.Item(89, 9).Value = 2844.3333  # ALC!I89
$ws.Cells.Item(89, 10).Value = 3495.8  # ALC!J89
$ws.Cells.Item(89, 11).Value = 14221.6665  # ALC!K89
$ws.Cells.Item(89, 12).Value = 17479  # ALC!L89
$ws.Cells.Item(89, 13).Value = -8605.6665  # ALC!M89
$ws.Cells.Item(89, 14).Value = -28711  # ALC!N89
$ws.Cells.Item(113, 8).Value = 5143  # ALC!H113
$ws.Cells.Item(113, 9).Value = 3495  # ALC!I113
$ws.Cells.Item(113, 11).Value = 3495  # ALC!K113
$ws.Cells.Item(113, 13).Value = -241  # ALC!M113
$ws.Cells.Item(118, 8).Value = 590  # ALC!H118
$ws.Cells.Item(118, 9).Value = 590  # ALC!I118
$ws.Cells.Item(118, 11).Value = 1770  # ALC!K118
$ws.Cells.Item(118, 13).Value = -113  # ALC!M118
$ws.Cells.Item(125, 8).Value = 785.2  # ALC!H125
$ws.Cells.Item(125, 9).Value = 681.5  # ALC!I125
$ws.Cells.Item(125, 11).Value = 6133.5  # ALC!K125
$ws.Cells.Item(125, 13).Value = -3673.5  # ALC!M125
$ws.Cells.Item(138, 8).Value = 2637.0356  # ALC!H138
$ws.Cells.Item(138, 10).Value = 3402  # ALC!J138
$ws.Cells.Item(138, 12).Value = 10206  # ALC!L138
$ws.Cells.Item(138, 14).Value = -20486  # ALC!N138
$ws.Cells.Item(141, 8).Value = 4071  # ALC!H141
$ws.Cells.Item(141, 9).Value = 4596.6665  # ALC!I141
$ws.Cells.Item(141, 11).Value = 13789.9995  # ALC!K141
$ws.Cells.Item(141, 13).Value = -8609.999500000002  # ALC!M141

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 4387.5415  # ARM!H2
$ws.Cells.Item(2, 9).Value = 3971.7144  # ARM!I2
$ws.Cells.Item(2, 10).Value = 4969.7  # ARM!J2
$ws.Cells.Item(2, 11).Value = 3971.7144  # ARM!K2
$ws.Cells.Item(2, 12).Value = 4969.7  # ARM!L2
$ws.Cells.Item(2, 13).Value = -3858.7144  # ARM!M2
$ws.Cells.Item(2, 14).Value = -5195.7  # ARM!N2
$ws.Cells.Item(74, 8).Value = 18702.738  # ARM!H74
$ws.Cells.Item(74, 9).Value = 22803.705  # ARM!I74
$ws.Cells.Item(74, 11).Value = 22803.705  # ARM!K74
$ws.Cells.Item(74, 13).Value = -21929.705  # ARM!M74
$ws.Cells.Item(77, 8).Value = 18702.738  # ARM!H77
$ws.Cells.Item(77, 9).Value = 22803.705  # ARM!I77
$ws.Cells.Item(77, 11).Value = 114018.525  # ARM!K77
$ws.Cells.Item(77, 13).Value = -109650.525  # ARM!M77
$ws.Cells.Item(116, 8).Value = 4387.5415  # ARM!H116
$ws.Cells.Item(116, 9).Value = 3971.7144  # ARM!I116
$ws.Cells.Item(116, 10).Value = 4969.7  # ARM!J116
$ws.Cells.Item(116, 11).Value = 3971.7144  # ARM!K116
$ws.Cells.Item(116, 12).Value = 4969.7  # ARM!L116
$ws.Cells.Item(116, 13).Value = -1677.7144  # ARM!M116
$ws.Cells.Item(116, 14).Value = -9557.7  # ARM!N116
$ws.Cells.Item(122, 8).Value = 1000  # ARM!H122
$ws.Cells.Item(122, 9).Value = 1000  # ARM!I122
$ws.Cells.Item(122, 10).Value = 0  # ARM!J122
$ws.Cells.Item(122, 11).Value = 3000  # ARM!K122
$ws.Cells.Item(122, 12).Value = 0  # ARM!L122
$ws.Cells.Item(122, 13).Value = -550  # ARM!M122
$ws.Cells.Item(122, 14).ClearContents()  # ARM!N122

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 4387.5415  # BSM!H3
$ws.Cells.Item(3, 9).Value = 3971.7144  # BSM!I3
$ws.Cells.Item(3, 10).Value = 4969.7  # BSM!J3
$ws.Cells.Item(3, 11).Value = 3971.7144  # BSM!K3
$ws.Cells.Item(3, 12).Value = 4969.7  # BSM!L3
$ws.Cells.Item(3, 13).Value = -3857.7144  # BSM!M3
$ws.Cells.Item(3, 14).Value = -5197.7  # BSM!N3
$ws.Cells.Item(20, 8).Value = 10226.294  # BSM!H20
$ws.Cells.Item(20, 9).Value = 13249.2  # BSM!I20
$ws.Cells.Item(20, 11).Value = 13249.2  # BSM!K20
$ws.Cells.Item(20, 13).Value = -13002.2  # BSM!M20
$ws.Cells.Item(86, 8).Value = 1253357.9  # BSM!H86
$ws.Cells.Item(86, 10).Value = 4143.8335  # BSM!J86
$ws.Cells.Item(86, 12).Value = 4143.8335  # BSM!L86
$ws.Cells.Item(86, 14).Value = -6389.8335  # BSM!N86
$ws.Cells.Item(89, 8).Value = 1253357.9  # BSM!H89
$ws.Cells.Item(89, 10).Value = 4143.8335  # BSM!J89
$ws.Cells.Item(89, 12).Value = 20719.1675  # BSM!L89
$ws.Cells.Item(89, 14).Value = -31951.1675  # BSM!N89
$ws.Cells.Item(105, 8).Value = 6181.091  # BSM!H105
$ws.Cells.Item(105, 9).Value = 6142.143  # BSM!I105
$ws.Cells.Item(105, 11).Value = 6142.143  # BSM!K105
$ws.Cells.Item(105, 13).Value = -4395.143  # BSM!M105
$ws.Cells.Item(107, 8).Value = 1905.4286  # BSM!H107
$ws.Cells.Item(107, 9).Value = 467.8  # BSM!I107
$ws.Cells.Item(107, 10).Value = 5499.5  # BSM!J107
$ws.Cells.Item(107, 11).Value = 467.8  # BSM!K107
$ws.Cells.Item(107, 12).Value = 5499.5  # BSM!L107
$ws.Cells.Item(107, 13).Value = 1452.2  # BSM!M107
$ws.Cells.Item(107, 14).Value = -9339.5  # BSM!N107

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(99, 8).Value = 10636.306  # CRP!H99
$ws.Cells.Item(99, 9).Value = 7316.3335  # CRP!I99
$ws.Cells.Item(99, 10).Value = 11300.3  # CRP!J99
$ws.Cells.Item(99, 11).Value = 7316.3335  # CRP!K99
$ws.Cells.Item(99, 12).Value = 11300.3  # CRP!L99
$ws.Cells.Item(99, 13).Value = -5818.3335  # CRP!M99
$ws.Cells.Item(99, 14).Value = -14296.3  # CRP!N99
$ws.Cells.Item(107, 8).Value = 1574.2858  # CRP!H107
$ws.Cells.Item(107, 9).Value = 1658.4  # CRP!I107
$ws.Cells.Item(107, 10).Value = 1364  # CRP!J107
$ws.Cells.Item(107, 11).Value = 1658.4  # CRP!K107
$ws.Cells.Item(107, 12).Value = 1364  # CRP!L107
$ws.Cells.Item(107, 13).Value = 261.5999999999999  # CRP!M107
$ws.Cells.Item(107, 14).Value = -5204  # CRP!N107
$ws.Cells.Item(126, 8).Value = 10636.306  # CRP!H126
$ws.Cells.Item(126, 9).Value = 7316.3335  # CRP!I126
$ws.Cells.Item(126, 10).Value = 11300.3  # CRP!J126
$ws.Cells.Item(126, 11).Value = 21949.0005  # CRP!K126
$ws.Cells.Item(126, 12).Value = 33900.89999999999  # CRP!L126
$ws.Cells.Item(126, 13).Value = -19479.0005  # CRP!M126
$ws.Cells.Item(126, 14).Value = -38840.89999999999  # CRP!N126

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(7, 8).Value = 111439  # CUL!H7
$ws.Cells.Item(7, 9).Value = 60  # CUL!I7
$ws.Cells.Item(7, 11).Value = 180  # CUL!K7
$ws.Cells.Item(7, 13).Value = -68  # CUL!M7
$ws.Cells.Item(131, 8).Value = 2599.2273  # CUL!H131
$ws.Cells.Item(131, 10).Value = 2831.0789  # CUL!J131
$ws.Cells.Item(131, 12).Value = 8493.2367  # CUL!L131
$ws.Cells.Item(131, 14).Value = -18573.2367  # CUL!N131

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 3690  # GSM!H122
$ws.Cells.Item(122, 9).Value = 3544.5557  # GSM!I122
$ws.Cells.Item(122, 11).Value = 10633.6671  # GSM!K122
$ws.Cells.Item(122, 13).Value = -8183.667099999999  # GSM!M122
$ws.Cells.Item(126, 8).Value = 5663.636  # GSM!H126
$ws.Cells.Item(126, 10).Value = 6069.2856  # GSM!J126
$ws.Cells.Item(126, 12).Value = 18207.8568  # GSM!L126
$ws.Cells.Item(126, 14).Value = -23147.8568  # GSM!N126

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 9249.571  # LTW!H7
$ws.Cells.Item(7, 9).Value = 7499.5  # LTW!I7
$ws.Cells.Item(7, 10).Value = 9541.25  # LTW!J7
$ws.Cells.Item(7, 11).Value = 7499.5  # LTW!K7
$ws.Cells.Item(7, 12).Value = 9541.25  # LTW!L7
$ws.Cells.Item(7, 13).Value = -7387.5  # LTW!M7
$ws.Cells.Item(7, 14).Value = -9765.25  # LTW!N7
$ws.Cells.Item(40, 8).Value = 7000  # LTW!H40
$ws.Cells.Item(40, 9).Value = 7000  # LTW!I40
$ws.Cells.Item(40, 10).Value = 0  # LTW!J40
$ws.Cells.Item(40, 11).Value = 7000  # LTW!K40
$ws.Cells.Item(40, 12).Value = 0  # LTW!L40
$ws.Cells.Item(40, 13).Value = -6864  # LTW!M40
$ws.Cells.Item(40, 14).ClearContents()  # LTW!N40
$ws.Cells.Item(54, 8).Value = 40000  # LTW!H54
$ws.Cells.Item(54, 10).Value = 40000  # LTW!J54
$ws.Cells.Item(54, 12).Value = 40000  # LTW!L54
$ws.Cells.Item(54, 14).Value = -41288  # LTW!N54
$ws.Cells.Item(82, 8).Value = 2341.5293  # LTW!H82
$ws.Cells.Item(82, 9).Value = 2019.1111  # LTW!I82
$ws.Cells.Item(82, 10).Value = 2704.25  # LTW!J82
$ws.Cells.Item(82, 11).Value = 2019.1111  # LTW!K82
$ws.Cells.Item(82, 12).Value = 2704.25  # LTW!L82
$ws.Cells.Item(82, 13).Value = -1658.1111  # LTW!M82
$ws.Cells.Item(82, 14).Value = -3426.25  # LTW!N82
$ws.Cells.Item(85, 8).Value = 2341.5293  # LTW!H85
$ws.Cells.Item(85, 9).Value = 2019.1111  # LTW!I85
$ws.Cells.Item(85, 10).Value = 2704.25  # LTW!J85
$ws.Cells.Item(85, 11).Value = 2019.1111  # LTW!K85
$ws.Cells.Item(85, 12).Value = 2704.25  # LTW!L85
$ws.Cells.Item(85, 13).Value = -771.1111000000001  # LTW!M85
$ws.Cells.Item(85, 14).Value = -5200.25  # LTW!N85
$ws.Cells.Item(122, 8).Value = 3438.318  # LTW!H122
$ws.Cells.Item(122, 9).Value = 4022.6667  # LTW!I122
$ws.Cells.Item(122, 10).Value = 3033.7693  # LTW!J122
$ws.Cells.Item(122, 11).Value = 12068.0001  # LTW!K122
$ws.Cells.Item(122, 12).Value = 9101.3079  # LTW!L122
$ws.Cells.Item(122, 13).Value = -9618.000100000001  # LTW!M122
$ws.Cells.Item(122, 14).Value = -14001.3079  # LTW!N122
$ws.Cells.Item(126, 8).Value = 9249.571  # LTW!H126
$ws.Cells.Item(126, 9).Value = 7499.5  # LTW!I126
$ws.Cells.Item(126, 10).Value = 9541.25  # LTW!J126
$ws.Cells.Item(126, 11).Value = 22498.5  # LTW!K126
$ws.Cells.Item(126, 12).Value = 28623.75  # LTW!L126
$ws.Cells.Item(126, 13).Value = -20028.5  # LTW!M126
$ws.Cells.Item(126, 14).Value = -33563.75  # LTW!N126
$ws.Cells.Item(132, 8).Value = 7862.5713  # LTW!H132
$ws.Cells.Item(132, 9).Value = 7098.8184  # LTW!I132
$ws.Cells.Item(132, 11).Value = 21296.4552  # LTW!K132
$ws.Cells.Item(132, 13).Value = -18766.4552  # LTW!M132

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(34, 8).Value = 64910  # WVR!H34
$ws.Cells.Item(34, 10).Value = 76131  # WVR!J34
$ws.Cells.Item(34, 12).Value = 76131  # WVR!L34
$ws.Cells.Item(34, 14).Value = -76537  # WVR!N34
$ws.Cells.Item(53, 8).Value = 11249  # WVR!H53
$ws.Cells.Item(53, 10).Value = 0  # WVR!J53
$ws.Cells.Item(53, 12).Value = 0  # WVR!L53
$ws.Cells.Item(53, 14).ClearContents()  # WVR!N53
$ws.Cells.Item(74, 8).Value = 4570.1665  # WVR!H74
$ws.Cells.Item(74, 10).Value = 2564  # WVR!J74
$ws.Cells.Item(74, 12).Value = 2564  # WVR!L74
$ws.Cells.Item(74, 14).Value = -4436  # WVR!N74
$ws.Cells.Item(77, 8).Value = 4570.1665  # WVR!H77
$ws.Cells.Item(77, 10).Value = 2564  # WVR!J77
$ws.Cells.Item(77, 12).Value = 7692  # WVR!L77
$ws.Cells.Item(77, 14).Value = -17052  # WVR!N77
$ws.Cells.Item(107, 8).Value = 5052098  # WVR!H107
$ws.Cells.Item(107, 9).Value = 1080.3334  # WVR!I107
$ws.Cells.Item(107, 10).Value = 15875707  # WVR!J107
$ws.Cells.Item(107, 11).Value = 3241.0002  # WVR!K107
$ws.Cells.Item(107, 12).Value = 47627121  # WVR!L107
$ws.Cells.Item(107, 13).Value = -1321.0002  # WVR!M107
$ws.Cells.Item(107, 14).Value = -47630961  # WVR!N107
$ws.Cells.Item(126, 8).Value = 7117.5757  # WVR!H126
$ws.Cells.Item(126, 9).Value = 5582.8076  # WVR!I126
$ws.Cells.Item(126, 10).Value = 12818.143  # WVR!J126
$ws.Cells.Item(126, 11).Value = 16748.4228  # WVR!K126
$ws.Cells.Item(126, 12).Value = 38454.429  # WVR!L126
$ws.Cells.Item(126, 13).Value = -14278.4228  # WVR!M126
$ws.Cells.Item(126, 14).Value = -43394.429  # WVR!N126
$ws.Cells.Item(131, 8).Value = 0  # WVR!H131
$ws.Cells.Item(131, 10).Value = 0  # WVR!J131
$ws.Cells.Item(131, 12).Value = 0  # WVR!L131
$ws.Cells.Item(131, 14).ClearContents()  # WVR!N131
$ws.Cells.Item(132, 8).Value = 153061.48  # WVR!H132
$ws.Cells.Item(132, 9).Value = 289827.4  # WVR!I132
$ws.Cells.Item(132, 11).Value = 869482.2000000001  # WVR!K132
$ws.Cells.Item(132, 13).Value = -866952.2000000001  # WVR!M132
